# Rename the two "Weekly" worksheets to their "Monthly" counterparts and
# update every cell that referenced the old weekly-oriented values, matching
# the authored change: this Scottish payroll input workbook was repurposed
# from the Week-4 suite to a Monthly-cadence run, and the "EMP 107" dummy
# automation guard employee record was swapped for "EMP 105".

$wb = $excel.ActiveWorkbook

$wsFirst   = $wb.Worksheets.Item("first")
$wsWeekly  = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsReports = $wb.Worksheets.Item("TestReports")

# --- Rename sheets (tab names + the "first" lookup sheet's text rows) ---
$wsWeekly.Name  = "GeneralTaxRateMonthly"
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# --- Swap the "DO NOT TOUCH AUTOMATION EMP 107" placeholder employee for
#     "DO NOT TOUCH AUTOMATION EMP 105" everywhere it appears ---
$wsWeekly.Range("A2").Value  = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# --- Restore the view/selection state recorded for each sheet ---
$wsFirst.Range("F5").Select()
$wsWeekly.Range("A4:XFD15").Select()
$wsProcess.Range("F9").Select()
$wsReports.Range("B10").Select()
